$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3").Value = 0.0082
$ws.Range("J7").Value = -0.0365
$ws.Range("K7").Value = -0.0354
$ws.Range("O7").Value = -0.0073
$ws.Range("Q7").Value = -0.0687
$ws.Range("R7").Value = -0.0299
$ws.Range("J9").Value = -0.1777
$ws.Range("K9").Value = -0.4306
$ws.Range("J11").Value = -0.0587
$ws.Range("K11").Value = -0.0852
$ws.Range("L11").Value = -0.0915
$ws.Range("M11").Value = -0.076
$ws.Range("N11").Value = -0.0726
$ws.Range("O11").Value = -0.0625
$ws.Range("P11").Value = -0.0492
$ws.Range("Q11").Value = -0.0529
$ws.Range("R11").Value = -0.0338
$ws.Range("J13").Value = -0.0735
$ws.Range("K13").Value = -0.0839
$ws.Range("L13").Value = 0.0636
$ws.Range("M13").Value = -0.0358
$ws.Range("N13").Value = -0.0213
$ws.Range("O13").Value = -0.0458
$ws.Range("P13").Value = -0.0706
$ws.Range("Q13").Value = 0.0349
$ws.Range("R13").Value = 0.0334
$ws.Range("J15").Value = -0.453
$ws.Range("K15").Value = -0.2229
$ws.Range("L15").Value = -0.0924
$ws.Range("M15").Value = -0.0679
$ws.Range("N15").Value = 0.0287
$ws.Range("O15").Value = 0.3272
$ws.Range("P15").Value = 0.321
$ws.Range("Q15").Value = 0.3382
$ws.Range("R15").Value = 0.3453
$ws.Range("J17").Value = -0.0907
$ws.Range("K17").Value = -0.0624
$ws.Range("L17").Value = -0.0388
$ws.Range("M17").Value = -0.0157
$ws.Range("N17").Value = -0.0178
$ws.Range("O17").Value = -0.0394
$ws.Range("P17").Value = -0.0344
$ws.Range("Q17").Value = -0.0284
$ws.Range("J19").Value = -0.0927
$ws.Range("K19").Value = -0.0806
$ws.Range("L19").Value = -0.0939
$ws.Range("M19").Value = -0.0927
$ws.Range("N19").Value = -0.0817
$ws.Range("J21").Value = 0.0071
$ws.Range("J23").Value = -0.5943
$ws.Range("K23").Value = -0.1432
$ws.Range("L23").Value = -0.0174
$ws.Range("M23").Value = 0.0456
$ws.Range("J27").Value = -0.0708
$ws.Range("K27").Value = -0.0589
$ws.Range("L27").Value = -0.5233
$ws.Range("M27").Value = -0.5717
$ws.Range("N27").Value = -0.2222
$ws.Range("O27").Value = -0.1758
$ws.Range("P27").Value = -0.3206
$ws.Range("Q27").Value = -0.2399
$ws.Range("R27").Value = -0.1032
$ws.Range("J29").Value = -0.3956
$ws.Range("K29").Value = -0.3902
$ws.Range("L29").Value = -0.2835
$ws.Range("M29").Value = -0.1285
$ws.Range("N29").Value = -0.0067
$ws.Range("O29").Value = -0.0009
$ws.Range("P29").Value = -0.0002
$ws.Range("Q29").Value = 0.0001
$ws.Range("R29").Value = 0.0001
$ws.Range("J31").Value = -2.1018
$ws.Range("K31").Value = -1.8326
$ws.Range("L31").Value = -2.5519
$ws.Range("M31").Value = -1.0643
$ws.Range("N31").Value = -0.4702
$ws.Range("O31").Value = -0.1682
$ws.Range("P31").Value = -0.368
$ws.Range("Q31").Value = -0.3384
$ws.Range("R31").Value = -0.0074
$ws.Range("J35").Value = -0.1803
$ws.Range("K35").Value = -0.4223
$ws.Range("J43").Value = 0.1361
$ws.Range("K43").Value = 0.3552
$ws.Range("J45").Value = -0.0064
$ws.Range("L45").Value = -0.0086
$ws.Range("M45").Value = -0.0001
$ws.Range("O45").Value = 0.0039
$ws.Range("Q45").Value = 0.0064
$ws.Range("J47").Value = -0.0786
$ws.Range("K47").Value = 0.0049
$ws.Range("L47").Value = -0.0816
$ws.Range("M47").Value = 0.1253
$ws.Range("N47").Value = 0.111
$ws.Range("O47").Value = 0.0968
$ws.Range("P47").Value = 0.083
$ws.Range("Q47").Value = -0.0598
$ws.Range("R47").Value = -0.0478
$ws.Range("J49").Value = 0.004
$ws.Range("K49").Value = -0.0043
$ws.Range("L49").Value = 0.0616
$ws.Range("M49").Value = 0.0007
$ws.Range("N49").Value = 0.0194
$ws.Range("O49").Value = 0.0075
$ws.Range("P49").Value = -0.0079
$ws.Range("Q49").Value = -0.0151
$ws.Range("K57").Value = -0.0403
$ws.Range("L57").Value = -0.0274
$ws.Range("M57").Value = -0.0163
$ws.Range("J59").Value = 0.0014
$ws.Range("J63").Value = -0.0001
$ws.Range("K63").Value = -0.0001
$ws.Range("Q63").Value = 0
$ws.Range("J65").Value = 0.038
$ws.Range("K65").Value = -0.0728
$ws.Range("J67").Value = -0.0029
$ws.Range("K67").Value = -0.0029
$ws.Range("L67").Value = -0.0029
$ws.Range("M67").Value = -0.003
$ws.Range("N67").Value = -0.003
$ws.Range("O67").Value = -0.003
$ws.Range("P67").Value = -0.0029
$ws.Range("Q67").Value = -0.0029
$ws.Range("R67").Value = -0.0029
$ws.Range("J69").Value = -0.0181
$ws.Range("K69").Value = -0.0198
$ws.Range("L69").Value = -0.0216
$ws.Range("M69").Value = -0.0208
$ws.Range("N69").Value = -0.0038
$ws.Range("O69").Value = -0.0017
$ws.Range("P69").Value = 0.0004
$ws.Range("Q69").Value = 0.0006
$ws.Range("R69").Value = 0.0008
$ws.Range("J71").Value = 0.0563
$ws.Range("K71").Value = 0.0294
$ws.Range("L71").Value = 0.007
$ws.Range("M71").Value = 0.0168
$ws.Range("N71").Value = 0.0164
$ws.Range("O71").Value = 0.0159
$ws.Range("P71").Value = 0.0157
$ws.Range("Q71").Value = 0.0105
$ws.Range("R71").Value = -0.0104
$ws.Range("J73").Value = 0.0111
$ws.Range("K73").Value = -0.0032
$ws.Range("L73").Value = -0.0001
$ws.Range("M73").Value = -0.0039
$ws.Range("N73").Value = -0.0001
$ws.Range("J75").Value = -0.0002
$ws.Range("K75").Value = -0.0001
$ws.Range("L75").Value = -0.0001
$ws.Range("M75").Value = -0.0001
$ws.Range("N75").Value = 0
$ws.Range("J77").Value = 0.0014
$ws.Range("J79").Value = -0.0142
$ws.Range("K79").Value = -0.014
$ws.Range("L79").Value = -0.0139
$ws.Range("M79").Value = -0.0085
$ws.Range("J83").Value = 0.0027
$ws.Range("K83").Value = -0.0005
$ws.Range("L83").Value = 0.0002
$ws.Range("M83").Value = 0.0003
$ws.Range("N83").Value = 0
$ws.Range("O83").Value = 0.0001
$ws.Range("P83").Value = 0.0002
$ws.Range("Q83").Value = 0.0001
$ws.Range("R83").Value = -0.0001
$ws.Range("J85").Value = 0.002
$ws.Range("K85").Value = 0.0017
$ws.Range("L85").Value = 0.0014
$ws.Range("M85").Value = 0.0012
$ws.Range("N85").Value = 0.0006
$ws.Range("O85").Value = 0.0003
$ws.Range("P85").Value = 0.0005
$ws.Range("Q85").Value = 0.0004
$ws.Range("R85").Value = 0.0001
$ws.Range("J87").Value = 0.0217
$ws.Range("K87").Value = -0.0059
$ws.Range("L87").Value = -0.029
$ws.Range("M87").Value = -0.0072
$ws.Range("N87").Value = 0.0262
$ws.Range("O87").Value = 0.028
$ws.Range("P87").Value = 0.0309
$ws.Range("Q87").Value = 0.013
$ws.Range("R87").Value = -0.0111
$ws.Range("J91").Value = 0.0735
$ws.Range("K91").Value = -0.0728
$ws.Range("J99").Value = -0.0735
$ws.Range("K99").Value = 0.0727
$ws.Range("J101").Value = 0.0001
$ws.Range("K101").Value = 0.0001
$ws.Range("L101").Value = 0.0001
$ws.Range("M101").Value = 0.0001
$ws.Range("N101").Value = 0.0001
$ws.Range("O101").Value = 0.0001
$ws.Range("P101").Value = 0.0001
$ws.Range("Q101").Value = 0.0001
$ws.Range("R101").Value = 0.0001
$ws.Range("J103").Value = -0.0175
$ws.Range("K103").Value = -0.0179
$ws.Range("L103").Value = -0.0181
$ws.Range("M103").Value = -0.0071
$ws.Range("N103").Value = -0.003
$ws.Range("O103").Value = -0.0021
$ws.Range("P103").Value = -0.0012
$ws.Range("Q103").Value = 0.0015
$ws.Range("R103").Value = 0.0013
$ws.Range("J105").Value = 0.0389
$ws.Range("K105").Value = 0.022
$ws.Range("L105").Value = 0.019
$ws.Range("M105").Value = 0.0188
$ws.Range("N105").Value = 0.0186
$ws.Range("O105").Value = 0.0184
$ws.Range("P105").Value = 0.0183
$ws.Range("Q105").Value = 0.0027
